# Wrap the three event-date strings (Events table, "Date" column) in
# backticks, e.g. "26th November 2025, Wednesday" -> "`26th November 2025, Wednesday`".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item("EventsTable")
$dateCells = $tbl.ListColumns.Item("Date").DataBodyRange

for ($i = 1; $i -le $dateCells.Rows.Count; $i++) {
    $cell = $dateCells.Cells.Item($i, 1)
    $cell.Value = "``" + $cell.Value2 + "``"
}
